$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 147207041.707713
$ws.Range("D2").Value = 71.10764899999999

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 22751356.077136
$ws.Range("D3").Value = 5.494966
$ws.Range("E3").Value = 0.00449

# Row 4 - Residuals
$ws.Range("B4").Value = 687306338.6878279
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -612.26325
$ws.Range("H5").Value = -1133.926428
$ws.Range("I5").Value = -90.600071
$ws.Range("J5").Value = 0.016603

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = -150.272569
$ws.Range("H6").Value = -698.195769
$ws.Range("I6").Value = 397.650632
$ws.Range("J6").Value = 0.79494

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 461.990681
$ws.Range("H7").Value = 49.98903
$ws.Range("I7").Value = 873.992332
$ws.Range("J7").Value = 0.023543
